$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction4")

$values = 0,1,2,4,5,6,7,8,9,10,11,12,13,14,16,17,18,19,20,21,22,23,24,26,27,29,30,31,33,34,36,37,38,39,40,41,42,43,44,47,49,50,52,53,54,56,57,58,59,60,63,69,70,71,73,74

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
